$wb = $excel.ActiveWorkbook

# ===== Sheet 展览 =====
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 196
$ws1.Cells.Item(4, 6).Value = 53
$ws1.Cells.Item(5, 6).Value = 1625
$ws1.Cells.Item(6, 6).Value = 3236
$ws1.Cells.Item(7, 6).Value = 778
$ws1.Cells.Item(8, 6).Value = 1991
$ws1.Cells.Item(9, 6).Value = 1911
$ws1.Cells.Item(10, 6).Value = 980
$ws1.Cells.Item(11, 6).Value = 344
$ws1.Cells.Item(13, 6).Value = 1583
$ws1.Cells.Item(14, 6).Value = 337
$ws1.Cells.Item(16, 6).Value = 62

# Insert new row at position 17 (pushes Look Look.. down to 18, etc.)
$ws1.Rows.Item(17).Insert()
$ws1.Cells.Item(16, 1).Copy($ws1.Cells.Item(17, 1))

# Fill content for rows 17..28
$ws1.Cells.Item(17, 1).Value = 16
$ws1.Cells.Item(17, 2).NumberFormat = "@"
$ws1.Cells.Item(17, 2).Value = "2024-05-26"
$ws1.Cells.Item(17, 3).Value = "广州·孤独摇滚only2.0"
$ws1.Cells.Item(17, 4).Value = "黄边三横路一街1号 设计殿堂"
$ws1.Cells.Item(17, 5).Value = "2024.05.26 10:00-05.26 17:00"
$ws1.Cells.Item(17, 6).Value = 17
$ws1.Cells.Item(17, 7).Value = 60
$ws1.Cells.Item(17, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85122"
$ws1.Cells.Item(17, 9).Value = "//i0.hdslb.com/bfs/openplatform/202405/N4JHQSfF1714988273293.png"

$ws1.Cells.Item(18, 1).Value = 17
$ws1.Cells.Item(18, 2).NumberFormat = "@"
$ws1.Cells.Item(18, 2).Value = "2024-06-01"
$ws1.Cells.Item(18, 3).Value = "广州·Look Look动漫嘉年华"
$ws1.Cells.Item(18, 4).Value = "东沙大道16号 健康方舟"
$ws1.Cells.Item(18, 5).Value = "2024.06.01 10:00-06.02 17:30"
$ws1.Cells.Item(18, 6).Value = 1394
$ws1.Cells.Item(18, 7).Value = 52.2
$ws1.Cells.Item(18, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82319"
$ws1.Cells.Item(18, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/QrBvxNAX1712126496119.jpeg"

$ws1.Cells.Item(19, 1).Value = 18
$ws1.Cells.Item(19, 2).NumberFormat = "@"
$ws1.Cells.Item(19, 2).Value = "2024-06-01"
$ws1.Cells.Item(19, 3).Value = "广州·WIO JUMPONLY3.0"
$ws1.Cells.Item(19, 4).Value = "黄边三横路一街1号 设计殿堂"
$ws1.Cells.Item(19, 5).Value = "2024.06.01 10:00-06.02 18:00"
$ws1.Cells.Item(19, 6).Value = 496
$ws1.Cells.Item(19, 7).Value = 70
$ws1.Cells.Item(19, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84722"
$ws1.Cells.Item(19, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/FhaZLO921713774163735.jpeg"

$ws1.Cells.Item(20, 1).Value = 19
$ws1.Cells.Item(20, 2).NumberFormat = "@"
$ws1.Cells.Item(20, 2).Value = "2024-06-01"
$ws1.Cells.Item(20, 3).Value = "广州·第五届AP动漫嘉年华"
$ws1.Cells.Item(20, 4).Value = "西环路1号 广州岭南会展中心"
$ws1.Cells.Item(20, 5).Value = "2024.06.01 10:00-06.01 17:00"
$ws1.Cells.Item(20, 6).Value = 601
$ws1.Cells.Item(20, 7).Value = 55
$ws1.Cells.Item(20, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83462"
$ws1.Cells.Item(20, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/ZR2jKMOg1711076939687.jpeg"

$ws1.Cells.Item(21, 1).Value = 20
$ws1.Cells.Item(21, 2).NumberFormat = "@"
$ws1.Cells.Item(21, 2).Value = "2024-06-08"
$ws1.Cells.Item(21, 3).Value = "广州·原神ONLY·旅行盛宴"
$ws1.Cells.Item(21, 4).Value = "桥头大街248号2层 格乐利雅GALLERIA艺术中心(海珠店)"
$ws1.Cells.Item(21, 5).Value = "2024.06.08 10:00-06.08 17:00"
$ws1.Cells.Item(21, 6).Value = 295
$ws1.Cells.Item(21, 7).Value = 60
$ws1.Cells.Item(21, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84662"
$ws1.Cells.Item(21, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/GO8aQp3d1713755965690.jpeg"

$ws1.Cells.Item(22, 1).Value = 21
$ws1.Cells.Item(22, 2).NumberFormat = "@"
$ws1.Cells.Item(22, 2).Value = "2024-06-08"
$ws1.Cells.Item(22, 3).Value = "广州·喵物语动漫游戏嘉年华"
$ws1.Cells.Item(22, 4).Value = "中新广州知识城凤桐直街12号 知识城国际会展中心"
$ws1.Cells.Item(22, 5).Value = "2024.06.08 10:00-06.10 16:00"
$ws1.Cells.Item(22, 6).Value = 10484
$ws1.Cells.Item(22, 7).Value = 78
$ws1.Cells.Item(22, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83856"
$ws1.Cells.Item(22, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/TjfGg7SU1711611802842.png"

$ws1.Cells.Item(23, 1).Value = 22
$ws1.Cells.Item(23, 2).NumberFormat = "@"
$ws1.Cells.Item(23, 2).Value = "2024-06-08"
$ws1.Cells.Item(23, 3).Value = "广州·珠三角 2024 COMIC WORLD次元世界动漫游戏嘉年华"
$ws1.Cells.Item(23, 4).Value = "南洲路139号 小洲云文化艺术创意园"
$ws1.Cells.Item(23, 5).Value = "2024.06.08 10:00-06.10 17:00"
$ws1.Cells.Item(23, 6).Value = 9673
$ws1.Cells.Item(23, 7).Value = 70
$ws1.Cells.Item(23, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85020"
$ws1.Cells.Item(23, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/6g0jnqBP1714146665737.jpeg"

$ws1.Cells.Item(24, 1).Value = 23
$ws1.Cells.Item(24, 2).NumberFormat = "@"
$ws1.Cells.Item(24, 2).Value = "2024-06-08"
$ws1.Cells.Item(24, 3).Value = "广州·运动番6.0-排球少年之宿命召集"
$ws1.Cells.Item(24, 4).Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws1.Cells.Item(24, 5).Value = "2024.06.08 10:00-06.08 17:00"
$ws1.Cells.Item(24, 6).Value = 837
$ws1.Cells.Item(24, 7).Value = 60
$ws1.Cells.Item(24, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83911"
$ws1.Cells.Item(24, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/i3Ngrbko1712129717623.jpeg"

$ws1.Cells.Item(25, 1).Value = 24
$ws1.Cells.Item(25, 2).NumberFormat = "@"
$ws1.Cells.Item(25, 2).Value = "2024-06-09"
$ws1.Cells.Item(25, 3).Value = "广州·AI动漫展5.0"
$ws1.Cells.Item(25, 4).Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws1.Cells.Item(25, 5).Value = "2024.06.09 10:00-06.09 17:00"
$ws1.Cells.Item(25, 6).Value = 639
$ws1.Cells.Item(25, 7).Value = 55
$ws1.Cells.Item(25, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83910"
$ws1.Cells.Item(25, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/DG65B7Bq1712130246181.jpeg"

$ws1.Cells.Item(26, 1).Value = 25
$ws1.Cells.Item(26, 2).NumberFormat = "@"
$ws1.Cells.Item(26, 2).Value = "2024-06-22"
$ws1.Cells.Item(26, 3).Value = "广州·622排球少年only"
$ws1.Cells.Item(26, 4).Value = "岭南购物城内 广州OMG网红街"
$ws1.Cells.Item(26, 5).Value = "2024.06.22 10:00-06.22 17:30"
$ws1.Cells.Item(26, 6).Value = 1810
$ws1.Cells.Item(26, 7).Value = 68
$ws1.Cells.Item(26, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82974"
$ws1.Cells.Item(26, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/WMlOXSZn1710748067155.jpeg"

$ws1.Cells.Item(27, 1).Value = 26
$ws1.Cells.Item(27, 2).NumberFormat = "@"
$ws1.Cells.Item(27, 2).Value = "2024-07-06"
$ws1.Cells.Item(27, 3).Value = "广州·重生之道only"
$ws1.Cells.Item(27, 4).Value = "同泰路颐和山庄 颐和大酒店"
$ws1.Cells.Item(27, 5).Value = "2024.07.06 10:30-07.06 16:30"
$ws1.Cells.Item(27, 6).Value = 134
$ws1.Cells.Item(27, 7).Value = 75
$ws1.Cells.Item(27, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84896"
$ws1.Cells.Item(27, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/aJpJGAEc1713699622756.png"

$ws1.Cells.Item(28, 1).Value = 27
$ws1.Cells.Item(28, 2).NumberFormat = "@"
$ws1.Cells.Item(28, 2).Value = "2024-07-14"
$ws1.Cells.Item(28, 3).Value = "广州·火影only"
$ws1.Cells.Item(28, 4).Value = "人和镇蚌湖清河大街168号 人和园"
$ws1.Cells.Item(28, 5).Value = "2024.07.14 09:30-07.14 17:30"
$ws1.Cells.Item(28, 6).Value = 384
$ws1.Cells.Item(28, 7).Value = 78
$ws1.Cells.Item(28, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84815"
$ws1.Cells.Item(28, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/QLOhW4Nr1714384036670.png"

# ===== Sheet 演出 =====
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(5, 6).Value = 113

# ===== Sheet 本地生活 =====
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 57

# ===== Sheet 全部类型 =====
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value = 57
$ws4.Cells.Item(4, 6).Value = 196
$ws4.Cells.Item(6, 6).Value = 53
$ws4.Cells.Item(7, 6).Value = 1625
$ws4.Cells.Item(8, 6).Value = 3236
$ws4.Cells.Item(9, 6).Value = 778
$ws4.Cells.Item(10, 6).Value = 1991
$ws4.Cells.Item(11, 6).Value = 1911
$ws4.Cells.Item(12, 6).Value = 980
$ws4.Cells.Item(13, 6).Value = 344
$ws4.Cells.Item(15, 6).Value = 1583
$ws4.Cells.Item(16, 6).Value = 337
$ws4.Cells.Item(19, 6).Value = 62

# Insert new row at position 21 (pushes Look Look.. down to 22, etc.)
$ws4.Rows.Item(21).Insert()
$ws4.Cells.Item(20, 1).Copy($ws4.Cells.Item(21, 1))

# Fill content for rows 21..35
$ws4.Cells.Item(21, 1).Value = 20
$ws4.Cells.Item(21, 2).NumberFormat = "@"
$ws4.Cells.Item(21, 2).Value = "2024-05-26"
$ws4.Cells.Item(21, 3).Value = "广州·孤独摇滚only2.0"
$ws4.Cells.Item(21, 4).Value = "黄边三横路一街1号 设计殿堂"
$ws4.Cells.Item(21, 5).Value = "2024.05.26 10:00-05.26 17:00"
$ws4.Cells.Item(21, 6).Value = 18
$ws4.Cells.Item(21, 7).Value = 60
$ws4.Cells.Item(21, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85122"
$ws4.Cells.Item(21, 9).Value = "//i0.hdslb.com/bfs/openplatform/202405/N4JHQSfF1714988273293.png"

$ws4.Cells.Item(22, 1).Value = 21
$ws4.Cells.Item(22, 2).NumberFormat = "@"
$ws4.Cells.Item(22, 2).Value = "2024-06-01"
$ws4.Cells.Item(22, 3).Value = "广州·Look Look动漫嘉年华"
$ws4.Cells.Item(22, 4).Value = "东沙大道16号 健康方舟"
$ws4.Cells.Item(22, 5).Value = "2024.06.01 10:00-06.02 17:30"
$ws4.Cells.Item(22, 6).Value = 1394
$ws4.Cells.Item(22, 7).Value = 52.2
$ws4.Cells.Item(22, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82319"
$ws4.Cells.Item(22, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/QrBvxNAX1712126496119.jpeg"

$ws4.Cells.Item(23, 1).Value = 22
$ws4.Cells.Item(23, 2).NumberFormat = "@"
$ws4.Cells.Item(23, 2).Value = "2024-06-01"
$ws4.Cells.Item(23, 3).Value = "广州·WIO JUMPONLY3.0"
$ws4.Cells.Item(23, 4).Value = "黄边三横路一街1号 设计殿堂"
$ws4.Cells.Item(23, 5).Value = "2024.06.01 10:00-06.02 18:00"
$ws4.Cells.Item(23, 6).Value = 496
$ws4.Cells.Item(23, 7).Value = 70
$ws4.Cells.Item(23, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84722"
$ws4.Cells.Item(23, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/FhaZLO921713774163735.jpeg"

$ws4.Cells.Item(24, 1).Value = 23
$ws4.Cells.Item(24, 2).NumberFormat = "@"
$ws4.Cells.Item(24, 2).Value = "2024-06-01"
$ws4.Cells.Item(24, 3).Value = "广州·第五届AP动漫嘉年华"
$ws4.Cells.Item(24, 4).Value = "西环路1号 广州岭南会展中心"
$ws4.Cells.Item(24, 5).Value = "2024.06.01 10:00-06.01 17:00"
$ws4.Cells.Item(24, 6).Value = 601
$ws4.Cells.Item(24, 7).Value = 55
$ws4.Cells.Item(24, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83462"
$ws4.Cells.Item(24, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/ZR2jKMOg1711076939687.jpeg"

$ws4.Cells.Item(25, 1).Value = 24
$ws4.Cells.Item(25, 2).NumberFormat = "@"
$ws4.Cells.Item(25, 2).Value = "2024-06-08"
$ws4.Cells.Item(25, 3).Value = "广州·原神ONLY·旅行盛宴"
$ws4.Cells.Item(25, 4).Value = "桥头大街248号2层 格乐利雅GALLERIA艺术中心(海珠店)"
$ws4.Cells.Item(25, 5).Value = "2024.06.08 10:00-06.08 17:00"
$ws4.Cells.Item(25, 6).Value = 295
$ws4.Cells.Item(25, 7).Value = 60
$ws4.Cells.Item(25, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84662"
$ws4.Cells.Item(25, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/GO8aQp3d1713755965690.jpeg"

$ws4.Cells.Item(26, 1).Value = 25
$ws4.Cells.Item(26, 2).NumberFormat = "@"
$ws4.Cells.Item(26, 2).Value = "2024-06-08"
$ws4.Cells.Item(26, 3).Value = "广州·喵物语动漫游戏嘉年华"
$ws4.Cells.Item(26, 4).Value = "中新广州知识城凤桐直街12号 知识城国际会展中心"
$ws4.Cells.Item(26, 5).Value = "2024.06.08 10:00-06.10 16:00"
$ws4.Cells.Item(26, 6).Value = 10484
$ws4.Cells.Item(26, 7).Value = 78
$ws4.Cells.Item(26, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83856"
$ws4.Cells.Item(26, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/TjfGg7SU1711611802842.png"

$ws4.Cells.Item(27, 1).Value = 26
$ws4.Cells.Item(27, 2).NumberFormat = "@"
$ws4.Cells.Item(27, 2).Value = "2024-06-08"
$ws4.Cells.Item(27, 3).Value = "广州·珠三角 2024 COMIC WORLD次元世界动漫游戏嘉年华"
$ws4.Cells.Item(27, 4).Value = "南洲路139号 小洲云文化艺术创意园"
$ws4.Cells.Item(27, 5).Value = "2024.06.08 10:00-06.10 17:00"
$ws4.Cells.Item(27, 6).Value = 9673
$ws4.Cells.Item(27, 7).Value = 70
$ws4.Cells.Item(27, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85020"
$ws4.Cells.Item(27, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/6g0jnqBP1714146665737.jpeg"

$ws4.Cells.Item(28, 1).Value = 27
$ws4.Cells.Item(28, 2).NumberFormat = "@"
$ws4.Cells.Item(28, 2).Value = "2024-06-08"
$ws4.Cells.Item(28, 3).Value = "广州·运动番6.0-排球少年之宿命召集"
$ws4.Cells.Item(28, 4).Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws4.Cells.Item(28, 5).Value = "2024.06.08 10:00-06.08 17:00"
$ws4.Cells.Item(28, 6).Value = 837
$ws4.Cells.Item(28, 7).Value = 60
$ws4.Cells.Item(28, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83911"
$ws4.Cells.Item(28, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/i3Ngrbko1712129717623.jpeg"

$ws4.Cells.Item(29, 1).Value = 28
$ws4.Cells.Item(29, 2).NumberFormat = "@"
$ws4.Cells.Item(29, 2).Value = "2024-06-09"
$ws4.Cells.Item(29, 3).Value = "广州·AI动漫展5.0"
$ws4.Cells.Item(29, 4).Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws4.Cells.Item(29, 5).Value = "2024.06.09 10:00-06.09 17:00"
$ws4.Cells.Item(29, 6).Value = 639
$ws4.Cells.Item(29, 7).Value = 55
$ws4.Cells.Item(29, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83910"
$ws4.Cells.Item(29, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/DG65B7Bq1712130246181.jpeg"

$ws4.Cells.Item(30, 1).Value = 29
$ws4.Cells.Item(30, 2).NumberFormat = "@"
$ws4.Cells.Item(30, 2).Value = "2024-06-22"
$ws4.Cells.Item(30, 3).Value = "广州·622排球少年only"
$ws4.Cells.Item(30, 4).Value = "岭南购物城内 广州OMG网红街"
$ws4.Cells.Item(30, 5).Value = "2024.06.22 10:00-06.22 17:30"
$ws4.Cells.Item(30, 6).Value = 1810
$ws4.Cells.Item(30, 7).Value = 68
$ws4.Cells.Item(30, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82974"
$ws4.Cells.Item(30, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/WMlOXSZn1710748067155.jpeg"

$ws4.Cells.Item(31, 1).Value = 30
$ws4.Cells.Item(31, 2).NumberFormat = "@"
$ws4.Cells.Item(31, 2).Value = "2024-06-28"
$ws4.Cells.Item(31, 3).Value = "广州·奥斯卡·罗曼耶卓（O叔）钢琴独奏音乐会"
$ws4.Cells.Item(31, 4).Value = "晴波路33号 广州星海音乐厅"
$ws4.Cells.Item(31, 5).Value = "2024.06.28 20:00-06.28 21:30"
$ws4.Cells.Item(31, 6).Value = 113
$ws4.Cells.Item(31, 7).Value = 180
$ws4.Cells.Item(31, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84545"
$ws4.Cells.Item(31, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/XK8EYxGv1712890578712.jpeg"

$ws4.Cells.Item(32, 1).Value = 31
$ws4.Cells.Item(32, 2).NumberFormat = "@"
$ws4.Cells.Item(32, 2).Value = "2024-06-29"
$ws4.Cells.Item(32, 3).Value = "广州·《海上钢琴师》经典电影作品大型交响音乐会"
$ws4.Cells.Item(32, 4).Value = "东风中路299号 广州中山纪念堂"
$ws4.Cells.Item(32, 5).Value = "2024.06.29 20:00-06.29 21:40"
$ws4.Cells.Item(32, 6).Value = 30
$ws4.Cells.Item(32, 7).Value = 75
$ws4.Cells.Item(32, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84162"
$ws4.Cells.Item(32, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/OnWieQKe1712742593534.jpeg"

$ws4.Cells.Item(33, 1).Value = 32
$ws4.Cells.Item(33, 2).NumberFormat = "@"
$ws4.Cells.Item(33, 2).Value = "2024-07-06"
$ws4.Cells.Item(33, 3).Value = "广州·重生之道only"
$ws4.Cells.Item(33, 4).Value = "同泰路颐和山庄 颐和大酒店"
$ws4.Cells.Item(33, 5).Value = "2024.07.06 10:30-07.06 16:30"
$ws4.Cells.Item(33, 6).Value = 134
$ws4.Cells.Item(33, 7).Value = 75
$ws4.Cells.Item(33, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84896"
$ws4.Cells.Item(33, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/aJpJGAEc1713699622756.png"

$ws4.Cells.Item(34, 1).Value = 33
$ws4.Cells.Item(34, 2).NumberFormat = "@"
$ws4.Cells.Item(34, 2).Value = "2024-07-14"
$ws4.Cells.Item(34, 3).Value = "广州·火影only"
$ws4.Cells.Item(34, 4).Value = "人和镇蚌湖清河大街168号 人和园"
$ws4.Cells.Item(34, 5).Value = "2024.07.14 09:30-07.14 17:30"
$ws4.Cells.Item(34, 6).Value = 384
$ws4.Cells.Item(34, 7).Value = 78
$ws4.Cells.Item(34, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84815"
$ws4.Cells.Item(34, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/QLOhW4Nr1714384036670.png"

$ws4.Cells.Item(35, 1).Value = 34
$ws4.Cells.Item(35, 2).NumberFormat = "@"
$ws4.Cells.Item(35, 2).Value = "2024-08-30"
$ws4.Cells.Item(35, 3).Value = "广州·孟京辉经典戏剧作品·黄湘丽主演《一个陌生女人的来信》"
$ws4.Cells.Item(35, 4).Value = "广州市越秀区人民北路696号 广州友谊剧院"
$ws4.Cells.Item(35, 5).Value = "2024.08.30 19:30-08.31 16:30"
$ws4.Cells.Item(35, 6).Value = 3
$ws4.Cells.Item(35, 7).Value = 100
$ws4.Cells.Item(35, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84570"
$ws4.Cells.Item(35, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/SscDFm1z1713177818070.jpeg"
